$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - KOREA AEROSPACE (047810.KS)
$ws.Range("D2").Value = 882000
$ws.Range("E2").Value = 34.4
$ws.Range("F2").Value = 3.52
$ws.Range("K2").Value = 54.7
$ws.Range("N2").Value = 54.83846622768671

# Row 3 - LIG Nex1 (079550.KS)
$ws.Range("D3").Value = 181400
$ws.Range("E3").Value = 35.8
$ws.Range("F3").Value = 3.24
$ws.Range("I3").Value = 63
$ws.Range("J3").Value = 63
$ws.Range("K3").Value = 50.7
$ws.Range("N3").Value = 54.83846622768671

# Row 4 - HANWHA AEROSPACE (012450.KS)
$ws.Range("D4").Value = 46950
$ws.Range("E4").Value = 24.1
$ws.Range("F4").Value = 1.62
$ws.Range("K4").Value = 46.5
$ws.Range("N4").Value = 54.83846622768671

# Row 5 - HYUNDAI ROTEM (064350.KS)
$ws.Range("D5").Value = 105400
$ws.Range("E5").Value = 38.8
$ws.Range("F5").Value = -3.21
$ws.Range("K5").Value = 44.7
$ws.Range("N5").Value = 54.83846622768671

# Row 6 - HANWHA SYSTEMS (272210.KS)
$ws.Range("D6").Value = 370000
$ws.Range("E6").Value = 28.2
$ws.Range("F6").Value = -3.39
$ws.Range("K6").Value = 37.9
$ws.Range("N6").Value = 54.83846622768671
